# Updates cryptos list with latest prices/volumes (GitHub Actions data refresh)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.419.45"
$ws.Range("E2").Value = "  +0.19%  "

# Row 3
$ws.Range("D3").Value = "1.850.01"
$ws.Range("E3").Value = "  +0.19%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6308"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "

# Row 7
$ws.Range("E7").Value = "  +0.07%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07706"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.19%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2942"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.47%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.43%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07750"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.63%  "

# Row 12
$ws.Range("D12").Value = "1.851.44"
$ws.Range("E12").Value = "  +0.11%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.027"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.71%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001084"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.29%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.09%  "

# Row 17
$ws.Range("D17").Value = "2.102.14"
$ws.Range("E17").Value = "  -0.18%  "

# Row 18
$ws.Range("E18").Value = "  +0.56%  "

# Row 19
$ws.Range("D19").Value = "29.447.36"
$ws.Range("E19").Value = "  +0.15%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.69%  "

# Row 21
$ws.Range("E21").Value = "  +0.30%  "

# Row 22
$ws.Range("E22").Value = "  +0.09%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.455"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.11%  "

# Row 24
$ws.Range("E24").Value = "  +0.05%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.23%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1390"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.49%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.364"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.10%  "

# Row 28
$ws.Range("E28").Value = "  +0.25%  "

# Row 29
$ws.Range("E29").Value = "  +0.56%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.312"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.79%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05734"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.04%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.110"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.26%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.054"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.80%  "

# Row 34
$ws.Range("E34").Value = "  +0.64%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.160"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.45%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7106"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.04%  "

# Row 37
$ws.Range("E37").Value = "  -0.32%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.779"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.01%  "

# Row 39 (now VeChain; was Maker)
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01801"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.65%  "

# Row 40 (now Maker; was VeChain)
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "1.229.74"
$ws.Range("E40").Value = "  -2.39%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.489"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.25%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9125"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.45%  "

# Row 43
$ws.Range("E43").Value = "  +0.11%  "

# Row 44
$ws.Range("D44").Value = "2.011.03"
$ws.Range("E44").Value = "  -0.20%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.24%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.27%  "

# Row 47
$ws.Range("E47").Value = "  +4.52%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.153"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.32%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4015"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.54%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.019"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.81%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.686"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.25%  "

